$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 158, shifting existing rows 158-204 down to 161-207
$ws.Rows.Item(158).Resize(3).Insert()

# Row 158
$ws.Range("A158").Value = 3
$ws.Range("B158").Value = 'Femacal de La Calera'
$ws.Range("C158").Value = 'Coquimbo'
$ws.Range("D158").Value = 45027
$ws.Range("E158").Value = 5
$ws.Range("F158").Value = 'Fruta'
$ws.Range("G158").Value = 100107
$ws.Range("H158").Value = 'Otros'
$ws.Range("I158").Value = 100107011
$ws.Range("J158").Value = 'Tuna'
$ws.Range("K158").Value = 'Sin especificar'
$ws.Range("L158").Value = 'Especial'
$ws.Range("M158").Value = 56
$ws.Range("N158").Value = 16000
$ws.Range("O158").Value = 16000
$ws.Range("P158").Value = 16000
$ws.Range("Q158").Value = '$/caja 16 kilos'
$ws.Range("R158").Value = 'Provincia de Los Andes'
$ws.Range("S158").Value = 1000
$ws.Range("T158").Value = 16

# Row 159
$ws.Range("A159").Value = 3
$ws.Range("B159").Value = 'Femacal de La Calera'
$ws.Range("C159").Value = 'Coquimbo'
$ws.Range("D159").Value = 45027
$ws.Range("E159").Value = 5
$ws.Range("F159").Value = 'Fruta'
$ws.Range("G159").Value = 100107
$ws.Range("H159").Value = 'Otros'
$ws.Range("I159").Value = 100107011
$ws.Range("J159").Value = 'Tuna'
$ws.Range("K159").Value = 'Sin especificar'
$ws.Range("L159").Value = 'Primera'
$ws.Range("M159").Value = 50
$ws.Range("N159").Value = 14000
$ws.Range("O159").Value = 14000
$ws.Range("P159").Value = 14000
$ws.Range("Q159").Value = '$/caja 16 kilos'
$ws.Range("R159").Value = 'Provincia de Los Andes'
$ws.Range("S159").Value = 875
$ws.Range("T159").Value = 16

# Row 160
$ws.Range("A160").Value = 3
$ws.Range("B160").Value = 'Femacal de La Calera'
$ws.Range("C160").Value = 'Coquimbo'
$ws.Range("D160").Value = 45027
$ws.Range("E160").Value = 5
$ws.Range("F160").Value = 'Fruta'
$ws.Range("G160").Value = 100107
$ws.Range("H160").Value = 'Otros'
$ws.Range("I160").Value = 100107011
$ws.Range("J160").Value = 'Tuna'
$ws.Range("K160").Value = 'Sin especificar'
$ws.Range("L160").Value = 'Segunda'
$ws.Range("M160").Value = 57
$ws.Range("N160").Value = 12000
$ws.Range("O160").Value = 12000
$ws.Range("P160").Value = 12000
$ws.Range("Q160").Value = '$/caja 16 kilos'
$ws.Range("R160").Value = 'Provincia de Los Andes'
$ws.Range("S160").Value = 750
$ws.Range("T160").Value = 16
